# Update Name of Algo
# Apply updated numeric results to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value  = -21.43680000000003
$ws.Range("D5").Value  = -8.402299999999995
$ws.Range("E7").Value  = 11.9956
$ws.Range("D9").Value  = -8.678700000000003
$ws.Range("D11").Value = -8.457700000000004
$ws.Range("E11").Value = 13.16369999999999
$ws.Range("A21").Value = -21.2817
$ws.Range("D21").Value = -7.983200000000008
$ws.Range("E21").Value = 13.17880000000001
$ws.Range("A23").Value = -21.47710000000002
$ws.Range("A25").Value = -22.38630000000003
